$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (sheet rows 6-12), matching columns A-J:
# A: id_response, B: start_date, C: end_date, D: status, E: progress,
# F: duration_s, G: finished, H: recorded_date, I: q1, J: q2

$data = @(
    @(5, 44523.81165509259, 44523.81331018519, "IP Address", 100, 142, $true, 44523.81332175926, "3mihar", "ebola %>% `npivot_longer(``289``:last_col(), names_to = ""day"", values_to = ""cases"") %>% `ndrop_na()"),
    @(6, 44523.81145833334, 44523.81542824074, "IP Address", 100, 342, $true, 44523.81542824074, "1pogus", $null),
    @(7, 44518.81680555556, 44518.85743055555, "Spam", 50, 3510, $false, 44523.87980324074, "0hacar", $null),
    @(8, 44518.81865740741, 44518.82510416667, "Spam", 50, 557, $false, 44523.87980324074, "0hacar", $null),
    @(9, 44523.80055555556, 44523.80107638889, "IP Address", 50, 44, $false, 44523.87981481482, "3mihar", $null),
    @(10, 44518.81951388888, 44518.82061342592, "IP Address", 50, 94, $false, 44523.87987268518, "8Musou", $null),
    @(11, 44523.09898148148, 44523.09974537038, "Spam", 50, 66, $false, 44523.87990740741, "2nesch", $null)
)

$r = 6
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    if ($row[9] -ne $null) {
        $ws.Cells.Item($r, 10).Value = $row[9]
    }
    $r++
}

# Undo the auto row-height bump that Excel applies to row 6 because of the
# embedded newlines in J6's text, so the row keeps the default height.
$ws.Rows.Item(6).AutoFit()

